$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# BDD syntax uses "Then" instead of "Assert" for the final step keyword.
# Update the cell that held the "Assert" label...
$ws.Range("A7").Value = "Then"

# ...and the matching conditional-formatting rule that highlights it.
$fcs = $ws.Range("A1:XFD1048576").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    if ($fc.Formula1 -eq '="Assert"') {
        $fc.Formula1 = '="Then"'
    }
}
